$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Education section tweaks
#    Para "Bachelors of Economics" / date / "Bachelors of Economics"
#      -> "Bachelors in Economics" / date / "UWI"
#    Para "Masters of Social Science" / date / "UWI"
#      -> "Masters in Social Science" / date / "UWI" (unchanged)
# ---------------------------------------------------------------------------

# Locate the two education paragraphs by their first line of text.
$bachelorsIndex = 0
$mastersIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($bachelorsIndex -eq 0 -and $t.StartsWith("Bachelors of Economics")) {
        $bachelorsIndex = $i
    }
    if ($mastersIndex -eq 0 -and $t.StartsWith("Masters of Social Science")) {
        $mastersIndex = $i
    }
}

$bachelorsPara = $d.Paragraphs.Item($bachelorsIndex)

# First occurrence in the paragraph -> "Bachelors in Economics"
$r1 = $bachelorsPara.Range
$r1.Find.Execute("Bachelors of Economics", $true, $false, $false, $false, $false, `
                  $true, 0, $false, "Bachelors in Economics", 1)

# Remaining occurrence (the trailing, third line) -> "UWI"
$r2 = $bachelorsPara.Range
$r2.Find.Execute("Bachelors of Economics", $true, $false, $false, $false, $false, `
                  $true, 0, $false, "UWI", 1)

$mastersPara = $d.Paragraphs.Item($mastersIndex)
$r3 = $mastersPara.Range
$r3.Find.Execute("Masters of Social Science", $true, $false, $false, $false, $false, `
                  $true, 0, $false, "Masters in Social Science", 1)

# ---------------------------------------------------------------------------
# 2. Skills section: add five bulleted "Label: Value" lines after the
#    "Skills" heading and before the "Projects" heading.
# ---------------------------------------------------------------------------

$skillsIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "Skills") {
        $skillsIndex = $i
    }
}

$items = @(
    @("Programming Language: ", "PHP"),
    @("Web Tech Framework: ", "Django,Flask"),
    @("Script: ", "JavaScript"),
    @("Database And Orm: ", "MySQL"),
    @("Version Control: ", "Git")
)

# Create all the (still empty) paragraphs first, right after "Skills", while
# the insertion point has no stray bold/italic formatting to inherit.
$anchor = $d.Paragraphs.Item($skillsIndex).Range
$anchor.Collapse(0)
for ($k = 0; $k -lt $items.Count; $k++) {
    $anchor.InsertParagraphAfter()
}

# Now fill each new paragraph with its styled runs.
for ($k = 0; $k -lt $items.Count; $k++) {
    $idx = $skillsIndex + 1 + $k
    $boldText = $items[$k][0]
    $italicText = $items[$k][1]

    $p = $d.Paragraphs.Item($idx)
    $p.Style = "List Bullet"

    $r = $p.Range
    $r.Collapse(1)
    $r.InsertAfter($boldText)
    $r.Bold = 1

    $r.Collapse(0)
    $r.InsertAfter($italicText)
    $r.Bold = 0
    $r.Italic = 1
}
